# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Column G (header "K") is recalculated for each outing row and rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G ("K"), per regenerated s_vals.
$kValues = @{
    2  = 2
    3  = 0
    4  = 1
    5  = 2
    6  = 1
    7  = 3
    8  = 0
    9  = 2
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 1
    16 = 1
    17 = 3
    18 = 1
    19 = 0
    21 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
